$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.421
$ws.Range("C6").Value = -11.843
$ws.Range("C7").Value = -12.793
$ws.Range("D7").Value = -7.673999999999999
$ws.Range("C8").Value = -12.839
$ws.Range("D11").Value = -7.337999999999999
$ws.Range("D12").Value = -7.587999999999999
$ws.Range("E12").Value = 17.491
$ws.Range("E13").Value = 16.366
$ws.Range("E14").Value = 16.854
$ws.Range("D15").Value = -8.392999999999999
$ws.Range("C16").Value = -12.808
$ws.Range("E16").Value = 16.706
$ws.Range("E19").Value = 16.638
$ws.Range("C20").Value = -12.241
$ws.Range("D20").Value = -7.994
$ws.Range("E20").Value = 16.696
$ws.Range("C21").Value = -12.428
$ws.Range("D21").Value = -8.196999999999999
$ws.Range("D22").Value = -7.536000000000001
$ws.Range("E22").Value = 16.998
$ws.Range("D23").Value = -7.997
$ws.Range("C28").Value = -12.848
$ws.Range("C29").Value = -12.438
$ws.Range("D29").Value = -7.416000000000001
$ws.Range("C30").Value = -12.559
$ws.Range("C32").Value = -12.471
$ws.Range("D34").Value = -7.936999999999999
$ws.Range("E36").Value = 16.801
$ws.Range("C40").Value = -12.119
$ws.Range("D42").Value = -8.18
$ws.Range("D43").Value = -7.972999999999999
$ws.Range("E43").Value = 17.028
$ws.Range("D44").Value = -7.580999999999999
$ws.Range("D45").Value = -7.517
$ws.Range("C46").Value = -13.322
$ws.Range("D46").Value = -8.139999999999999
$ws.Range("E46").Value = 16.889
$ws.Range("D50").Value = -8.312999999999999
$ws.Range("E50").Value = 16.688
$ws.Range("C51").Value = -12.199
$ws.Range("D51").Value = -7.598999999999999
$ws.Range("C52").Value = -11.614
$ws.Range("C57").Value = -13.329
$ws.Range("D57").Value = -8.228999999999999
$ws.Range("C59").Value = -12.244
$ws.Range("C62").Value = -13.737
$ws.Range("D65").Value = -7.556999999999999
$ws.Range("C66").Value = -11.36
$ws.Range("D66").Value = -7.683000000000002
$ws.Range("D67").Value = -6.82
$ws.Range("C73").Value = -12.52
$ws.Range("C74").Value = -11.866
$ws.Range("E76").Value = 16.995
$ws.Range("C77").Value = -12.54
$ws.Range("D79").Value = -7.728
$ws.Range("D84").Value = -8.453000000000001
$ws.Range("D87").Value = -8.068000000000001
$ws.Range("C92").Value = -11.604
$ws.Range("D92").Value = -6.728
$ws.Range("E95").Value = 17.52
$ws.Range("D97").Value = -8.486000000000001
$ws.Range("E97").Value = 16.422
$ws.Range("E99").Value = 16.786
$ws.Range("C100").Value = -12.623
